$d = $word.ActiveDocument

# Remove all paragraphs except the last one; we will replace the last
# paragraph (and implicitly insert everything before it) via InsertXML.
$lastIndex = $d.Paragraphs.Count
if ($lastIndex -gt 1) {
    $delRange = $d.Range(0, $d.Paragraphs($lastIndex - 1).Range.End)
    $delRange.Delete()
}

$xml = @'
<pkg:xmlData xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Chris Cargile</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:left="720" w:hanging="720"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>2.2)</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t>Specification Languages (SL) describes a system diagrammat</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ically, in terms of: </w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:left="720" w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">- Nodes =&gt; </w:t>
      </w:r>
      <w:r>
        <w:t>functions (nodes B</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>,C</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> might nest in ‘A’, for example);   </w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:left="720" w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t>- Arcs =&gt; data flows; and</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>c</w:t>
      </w:r>
      <w:r>
        <w:t>hecks</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> for self-consistency might include: </w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:left="1440" w:hanging="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">-&gt; </w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">do the functions (which ===data transformation) always maintain the same invariant(s) and when nesting is illustrated, do the nested transformations </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>facilititate</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">/permit adhering to invariant of </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>nestee</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>?  (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ie</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>functionA</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>int</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> a , </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>int</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> b) { return (</w:t>
      </w:r>
      <w:r>
        <w:t>&lt;T&gt;d</w:t>
      </w:r>
      <w:r>
        <w:t>ouble)</w:t>
      </w:r>
      <w:r>
        <w:t>?</w:t>
      </w:r>
      <w:r>
        <w:t>--</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>functionB</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> ( a , b ) { return double } }</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>-&gt;</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t>Are null flows allowed/disallowed consistently?</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:left="1440" w:hanging="720"/>
      </w:pPr>
      <w:r>
        <w:t>-&gt;</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t>Can functions requiring/forbidding data/null receive their counterparts consistently OR shows ‘throwing a flag’ consistently?</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">-&gt; </w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Are recursively written functions portrayed the same, </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>always</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
<w:p/>
<w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="0"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>2.4)</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Comparing relative approaches of (4) methods to ensuring </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>serializability</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> based on (3) axes of verification techniques (A) – Pessimistic Inaccuracy;  (B) – Optimistic Inaccuracy;  (C) Simplified Properties…..</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="0"/>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
<w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="0"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t>(1): A</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>,B,C</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> are the same essentially as A,B are ==0 (there is no inaccuracy) due to C==true</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
      </w:pPr>
      <w:r>
        <w:t>(2) A is the case here since the analysis may reject some valid sets.  B is not the case as the analysis shou</w:t>
      </w:r>
      <w:r>
        <w:t>ld not give false positives.  Automated static analysis cannot prove logical correctness therefore C is false.</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
      </w:pPr>
      <w:r>
        <w:t>(3) A,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>B should be false</w:t>
      </w:r>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> as conformance to a specification would seem a simplified-property check (C==true); however, if the </w:t>
      </w:r>
    </w:p>
</pkg:xmlData>
'@

$d.Paragraphs(1).Range.InsertXML($xml)

Write-Output "paragraphs now: $($d.Paragraphs.Count)"
